$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 125000984
$ws.Range("I32").Value = 142857980
$ws.Range("J32").Value = 2000
$ws.Range("K32").Value = 142857980
$ws.Range("L32").Value = 2000
$ws.Range("M32").Value = -142857654
$ws.Range("N32").Value = -2652
# Row 40
$ws.Range("H40").Value = 1683.72
$ws.Range("I40").Value = 1557.5264
$ws.Range("J40").Value = 2083.3333
$ws.Range("K40").Value = 1557.5264
$ws.Range("L40").Value = 2083.3333
$ws.Range("M40").Value = -1382.5264
$ws.Range("N40").Value = -2433.3333
# Row 51
$ws.Range("H51").Value = 7324.9
$ws.Range("I51").Value = 11719.1
$ws.Range("J51").Value = 2930.7
$ws.Range("K51").Value = 11719.1
$ws.Range("L51").Value = 2930.7
$ws.Range("M51").Value = -11235.1
$ws.Range("N51").Value = -3898.7
# Row 74
$ws.Range("H74").Value = 5833
$ws.Range("I74").Value = 4750
$ws.Range("J74").Value = 7999
$ws.Range("K74").Value = 4750
$ws.Range("L74").Value = 7999
$ws.Range("M74").Value = -3814
$ws.Range("N74").Value = -9871
# Row 77
$ws.Range("H77").Value = 5833
$ws.Range("I77").Value = 4750
$ws.Range("J77").Value = 7999
$ws.Range("K77").Value = 23750
$ws.Range("L77").Value = 39995
$ws.Range("M77").Value = -19070
$ws.Range("N77").Value = -49355
# Row 99
$ws.Range("H99").Value = 708.7143
$ws.Range("I99").Value = 396
$ws.Range("J99").Value = 1490.5
$ws.Range("K99").Value = 1188
$ws.Range("L99").Value = 4471.5
$ws.Range("M99").Value = 310
$ws.Range("N99").Value = -7467.5
# Row 100
$ws.Range("H100").Value = 1571.7142
$ws.Range("I100").Value = 1528.2222
$ws.Range("J100").Value = 1650
$ws.Range("K100").Value = 1528.2222
$ws.Range("L100").Value = 1650
$ws.Range("M100").Value = -987.2221999999999
$ws.Range("N100").Value = -2732
# Row 137
$ws.Range("H137").Value = 1354.463
$ws.Range("I137").Value = 1510.3077
$ws.Range("J137").Value = 1305.0488
$ws.Range("K137").Value = 4530.9231
$ws.Range("L137").Value = 3915.1464
$ws.Range("M137").Value = -1980.9231
$ws.Range("N137").Value = -9015.1464
# Row 138
$ws.Range("H138").Value = 4500.2295
$ws.Range("I138").Value = 2704.1177
$ws.Range("J138").Value = 5194.1816
$ws.Range("K138").Value = 8112.353099999999
$ws.Range("L138").Value = 15582.5448
$ws.Range("M138").Value = -2972.353099999999
$ws.Range("N138").Value = -25862.5448
# Row 141
$ws.Range("H141").Value = 2700.6875
$ws.Range("I141").Value = 2174.1155
$ws.Range("J141").Value = 4982.5
$ws.Range("K141").Value = 6522.3465
$ws.Range("L141").Value = 14947.5
$ws.Range("M141").Value = -1342.3465
$ws.Range("N141").Value = -25307.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 23426.54
$ws.Range("I32").Value = 4433.431
$ws.Range("J32").Value = 78506.55
$ws.Range("K32").Value = 4433.431
$ws.Range("L32").Value = 78506.55
$ws.Range("M32").Value = -4146.431
$ws.Range("N32").Value = -79080.55
# Row 63
$ws.Range("H63").Value = 2510.8235
$ws.Range("I63").Value = 2054.1538
$ws.Range("J63").Value = 3995
$ws.Range("K63").Value = 2054.1538
$ws.Range("L63").Value = 3995
$ws.Range("M63").Value = -1368.1538
$ws.Range("N63").Value = -5367
# Row 66
$ws.Range("H66").Value = 2510.8235
$ws.Range("I66").Value = 2054.1538
$ws.Range("J66").Value = 3995
$ws.Range("K66").Value = 10270.769
$ws.Range("L66").Value = 19975
$ws.Range("M66").Value = -6838.769
$ws.Range("N66").Value = -26839
# Row 97
$ws.Range("H97").Value = 43662.832
$ws.Range("I97").Value = 67619.8
$ws.Range("J97").Value = 3734.5557
$ws.Range("K97").Value = 67619.8
$ws.Range("L97").Value = 3734.5557
$ws.Range("M97").Value = -67123.8
$ws.Range("N97").Value = -4726.5557
# Row 119
$ws.Range("H119").Value = 30248
$ws.Range("I119").Value = 10000
$ws.Range("J119").Value = 33140.57
$ws.Range("K119").Value = 10000
$ws.Range("L119").Value = 33140.57
$ws.Range("M119").Value = -5162
$ws.Range("N119").Value = -42816.57
# Row 132
$ws.Range("H132").Value = 2808.75
$ws.Range("I132").Value = 3686.9666
$ws.Range("J132").Value = 1611.1818
$ws.Range("K132").Value = 11060.8998
$ws.Range("L132").Value = 4833.5454
$ws.Range("M132").Value = -8530.899800000001
$ws.Range("N132").Value = -9893.545399999999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 39
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = None
$ws.Range("N39").ClearContents()
# Row 94
$ws.Range("H94").Value = 111402.445
$ws.Range("I94").Value = 111402.445
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 111402.445
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -110951.445
# Row 134
$ws.Range("H134").Value = 1669.3334
$ws.Range("I134").Value = 1624.8889
$ws.Range("J134").Value = 1936
$ws.Range("K134").Value = 4874.6667
$ws.Range("L134").Value = 5808
$ws.Range("M134").Value = -2339.6667
$ws.Range("N134").Value = -10878

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 15111.105
$ws.Range("I31").Value = 33359.16
$ws.Range("J31").Value = 2540.2222
$ws.Range("K31").Value = 33359.16
$ws.Range("L31").Value = 2540.2222
$ws.Range("M31").Value = -33064.16
$ws.Range("N31").Value = -3130.2222
# Row 33
$ws.Range("H33").Value = 5015.5
$ws.Range("I33").Value = 5015.5
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 5015.5
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = None
$ws.Range("N33").ClearContents()
# Row 34
$ws.Range("H34").Value = 15111.105
$ws.Range("I34").Value = 33359.16
$ws.Range("J34").Value = 2540.2222
$ws.Range("K34").Value = 33359.16
$ws.Range("L34").Value = 2540.2222
$ws.Range("M34").Value = -33157.16
$ws.Range("N34").Value = -2944.2222
# Row 58
$ws.Range("H58").Value = 13927.429
$ws.Range("I58").Value = 1797.7646
$ws.Range("J58").Value = 65478.5
$ws.Range("K58").Value = 1797.7646
$ws.Range("L58").Value = 65478.5
$ws.Range("M58").Value = -1594.7646
$ws.Range("N58").Value = -65884.5
# Row 132
$ws.Range("H132").Value = 2853.1333
$ws.Range("I132").Value = 2895.2727
$ws.Range("J132").Value = 2737.25
$ws.Range("K132").Value = 8685.8181
$ws.Range("L132").Value = 8211.75
$ws.Range("M132").Value = -6155.8181
$ws.Range("N132").Value = -13271.75
# Row 134
$ws.Range("H134").Value = 1769.6
$ws.Range("I134").Value = 1324.8
$ws.Range("J134").Value = 2659.2
$ws.Range("K134").Value = 3974.4
$ws.Range("L134").Value = 7977.599999999999
$ws.Range("M134").Value = -1439.4
$ws.Range("N134").Value = -13047.6
# Row 136
$ws.Range("H136").Value = 13927.429
$ws.Range("I136").Value = 1797.7646
$ws.Range("J136").Value = 65478.5
$ws.Range("K136").Value = 5393.293799999999
$ws.Range("L136").Value = 196435.5
$ws.Range("M136").Value = -2843.293799999999
$ws.Range("N136").Value = -201535.5

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 1965.8933
$ws.Range("I68").Value = 1269.4412
$ws.Range("J68").Value = 2543.439
$ws.Range("K68").Value = 3808.3236
$ws.Range("L68").Value = 7630.316999999999
$ws.Range("M68").Value = -2997.3236
$ws.Range("N68").Value = -9252.316999999999
# Row 71
$ws.Range("H71").Value = 1965.8933
$ws.Range("I71").Value = 1269.4412
$ws.Range("J71").Value = 2543.439
$ws.Range("K71").Value = 11424.9708
$ws.Range("L71").Value = 22890.951
$ws.Range("M71").Value = -7368.970799999999
$ws.Range("N71").Value = -31002.951

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 66668692
$ws.Range("I97").Value = 111113496
$ws.Range("J97").Value = 1486.8334
$ws.Range("K97").Value = 111113496
$ws.Range("L97").Value = 1486.8334
$ws.Range("M97").Value = -111113000
$ws.Range("N97").Value = -2478.8334
# Row 107
$ws.Range("H107").Value = 594500.9399999999
$ws.Range("I107").Value = 385.5
$ws.Range("J107").Value = 1122603.5
$ws.Range("K107").Value = 385.5
$ws.Range("L107").Value = 1122603.5
$ws.Range("M107").Value = 1534.5
$ws.Range("N107").Value = -1126443.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 1835.8667
$ws.Range("I82").Value = 1868.8889
$ws.Range("J82").Value = 1786.3334
$ws.Range("K82").Value = 1868.8889
$ws.Range("L82").Value = 1786.3334
$ws.Range("M82").Value = -1507.8889
$ws.Range("N82").Value = -2508.3334
# Row 85
$ws.Range("H85").Value = 1835.8667
$ws.Range("I85").Value = 1868.8889
$ws.Range("J85").Value = 1786.3334
$ws.Range("K85").Value = 1868.8889
$ws.Range("L85").Value = 1786.3334
$ws.Range("M85").Value = -620.8888999999999
$ws.Range("N85").Value = -4282.3334
# Row 119
$ws.Range("H119").Value = 31221.25
$ws.Range("I119").Value = 20000
$ws.Range("J119").Value = 37954
$ws.Range("K119").Value = 20000
$ws.Range("L119").Value = 37954
$ws.Range("M119").Value = -15162
$ws.Range("N119").Value = -47630
# Row 136
$ws.Range("H136").Value = 1919.0385
$ws.Range("I136").Value = 1495
$ws.Range("J136").Value = 4251.25
$ws.Range("K136").Value = 4485
$ws.Range("L136").Value = 12753.75
$ws.Range("M136").Value = -1935
$ws.Range("N136").Value = -17853.75

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 250002500
$ws.Range("I96").Value = 250002500
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 250002500
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -250001127
# Row 119
$ws.Range("H119").Value = 34842.5
$ws.Range("I119").Value = 20000
$ws.Range("J119").Value = 39790
$ws.Range("K119").Value = 20000
$ws.Range("L119").Value = 39790
$ws.Range("M119").Value = -15162
$ws.Range("N119").Value = -49466
# Row 136
$ws.Range("H136").Value = 823.53125
$ws.Range("I136").Value = 497.52173
$ws.Range("J136").Value = 1656.6666
$ws.Range("K136").Value = 1492.56519
$ws.Range("L136").Value = 4969.9998
$ws.Range("M136").Value = 1057.52173
$ws.Range("N136").Value = -10069.9998
